$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 to make room for the new first forecast entry.
$ws.Rows.Item(2).Insert()

# The inserted row picks up formatting from the row above (the header row).
# Re-apply the plain data-row formatting (as used by the rest of the table)
# by copying the format from the row that used to be row 2 (now row 3).
$ws.Range("A3:E3").Copy()
$ws.Range("A2:E2").PasteSpecial(-4122)

# Full dataset (A=date serial, B=y_0, C=y_0_forecast, D=y_1, E=y_1_forecast)
$data = @(
    @(39400, 2007, 5.896808312953783, 2008, 10.70171490310616),
    @(39765, 2008, 7.441962824572235, 2009, 10.14224506046018),
    @(40130, 2009, 6.277541464866987, 2010, 9.926356894615008),
    @(40494, 2010, 6.535114773304773, 2011, 6.670100100023779),
    @(40862, 2011, 5.12051970717502, 2012, 6.878526556821973),
    @(41228, 2012, 3.65682115264816, 2013, 4.783388268072009),
    @(41592, 2013, 2.943878639034381, 2014, 5.440628060653574),
    @(41957, 2014, 1.172679597477866, 2015, 0.7846309372076199),
    @(42321, 2015, 2.961845079861303, 2016, 2.521390070591267),
    @(42689, 2016, 2.508469427909898, 2017, 2.540874511056623),
    @(43053, 2017, 3.523703831572056, 2018, 3.025650759929999),
    @(43418, 2018, 1.178605266817589, 2019, 0.6761281928490348),
    @(43783, 2019, 3.047037961814492, 2020, 3.638425237116749),
    @(44159, 2020, -0.2228847697281378, 2021, 2.646620531497335),
    @(44525, 2021, -1.165854108406617, 2022, 4.306257314632855),
    @(44890, 2022, 2.501311189006916, 2023, 3.687525625406263),
    @(45254, 2023, 0.6753076481029074, 2024, 0.5720316833079497),
    @(45618, 2024, 2.039329803030121, 2025, 2.838786065810939)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $row++
}
